# "first crack at 3-4"
#
# The old plan had two separate rows for lesson 3-4 (rows 23 & 24):
#   row 23: "3-4a" / "BLE CapSense Remote control of LED Service 3-2a"
#   row 24: "3-4b" / "Add CapSense to remote control"
# These get collapsed into a single "3-4" lesson row. We do this by
# deleting row 23 (which shifts row 24, with its already-correct
# "Add CapSense to remote control" title, up into row 23) and then
# renaming that surviving row's lesson number and project fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "3-4a" row entirely; everything below (old row 24
# "3-4b" onward) shifts up by one.
$ws.Rows(23).Delete()

# The row that slid up into position 23 was "3-4b" / "Add CapSense to
# remote control" — rename it to the merged lesson "3-4" and give it
# its project name. Title (B), Script flag (E) already hold the right
# values ("Add CapSense to remote control" / "n").
$ws.Range("A23").Value = "3-4"
$ws.Range("D23").Value = "rc3-4 CapSense"

# Match the author's final selection in the sheet.
[void]$ws.Range("A23").Select()
